$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The feed that backs this sheet picked up a new case (6465) that slots in
# ahead of the existing -604 row (by report/claim date), plus a brand-new
# case (7359) at the end. Insert a fresh row 37 - this pushes the current
# row 37 (-604 / Gurruchaga 2126 / ...) down to row 38 untouched - then
# populate the new row 37 and append the new row 39.

$ws.Rows.Item(37).Insert()

function Set-CaseRow($ws, $r, $vals) {
    # Caso (A), Comuna (D) and OT (E) are numeric-looking ids and F. De
    # Reclamo (B) is a date-looking string, but this sheet keeps every one
    # of those as plain text - pre-format as Text so Excel doesn't
    # reinterpret them as a number/date. The remaining text columns are
    # already alphabetic, so a plain value assignment already lands as
    # text; Attachments (I), Coordenada_X/Y (M/N) are the genuinely
    # numeric columns.
    $textNumericCols = @(1,2,4,5)
    foreach ($c in $textNumericCols) {
        $ws.Cells.Item($r, $c).NumberFormat = "@"
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }

    $plainTextCols = @(3,6,7,8,10,11,12,15,16,17,18)
    foreach ($c in $plainTextCols) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }

    $ws.Cells.Item($r, 9).Value = $vals[8]
    $ws.Cells.Item($r, 13).Value = $vals[12]
    $ws.Cells.Item($r, 14).Value = $vals[13]
}

$row37 = @("6465", "8/28/2025", "AGUIRRE 368", "15", "809268249", "AYKO", "Pendiente", "Picada", 1, "Cambio", "Sin equipos", "Pasante", -58.434651, -34.598814, "Palermo", "Capital Sur", "CLI-O", "Fuera de Poligono OVL")
Set-CaseRow $ws 37 $row37

$row39 = @("7359", "9/29/2025", "VEGA, NICETO, CNEL. 4678", "14", "810056579", "AYKO", "Pendiente", "Columna inclinada con base corroida", 1, "Cambio", "Sin equipos", "Pasante", -58.430056, -34.593188, "Palermo", "Capital Sur", "VCR-F", "Fuera de Poligono OVL")
Set-CaseRow $ws 39 $row39
